$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.074.13"
$ws.Range("E2").Value = "  -1.81%  "
$ws.Range("D3").Value = "2.471.31"
$ws.Range("E3").Value = "  -2.16%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "519.40"
$ws.Range("E5").Value = "  -3.28%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.81"
$ws.Range("E6").Value = "  -3.82%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  -2.06%  "
$ws.Range("E9").Value = "  -2.10%  "
$ws.Range("E10").Value = "  -0.38%  "
$ws.Range("E11").Value = "  +0.14%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.343"
$ws.Range("D13").Value = "2.908.65"
$ws.Range("E13").Value = "  -2.09%  "
$ws.Range("D14").Value = "57.976.55"
$ws.Range("E14").Value = "  -1.85%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "22.32"
$ws.Range("E15").Value = "  -3.15%  "
$ws.Range("E16").Value = "  -2.07%  "
$ws.Range("D17").Value = "2.469.12"
$ws.Range("E17").Value = "  -2.32%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.77"
$ws.Range("E18").Value = "  -3.29%  "
$ws.Range("E19").Value = "  -2.35%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "320.56"
$ws.Range("E20").Value = "  -0.89%  "
$ws.Range("E21").Value = "  +0.14%  "
$ws.Range("E22").Value = "  -3.35%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "64.09"
$ws.Range("E23").Value = "  -2.67%  "
$ws.Range("E24").Value = "  -2.40%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.998"
$ws.Range("E25").Value = "  -0.17%  "
$ws.Range("E26").Value = "  -3.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.35"
$ws.Range("E27").Value = "  -2.60%  "
$ws.Range("D28").Value = "0.0₃0751"
$ws.Range("E28").Value = "  -3.01%  "
$ws.Range("B29").Value = "PancakeSwap"
$ws.Range("C29").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.69"
$ws.Range("E29").Value = "  -4.08%  "
$ws.Range("B30").Value = "Aptos"
$ws.Range("C30").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.30"
$ws.Range("E30").Value = "  -6.23%  "
$ws.Range("B31").Value = "Monero"
$ws.Range("C31").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "165.80"
$ws.Range("E31").Value = "  -1.47%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.16"
$ws.Range("E32").Value = "  -1.33%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("E33").Value = "  -0.01%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.997"
$ws.Range("E34").Value = "  -0.21%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "18.03"
$ws.Range("E35").Value = "  -2.14%  "
$ws.Range("E36").Value = "  -10.12%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.99"
$ws.Range("E37").Value = "  -3.09%  "
$ws.Range("E38").Value = "  -4.39%  "
$ws.Range("E39").Value = "  -3.30%  "
$ws.Range("E40").Value = "  -4.51%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "273.03"
$ws.Range("E41").Value = "  -4.09%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.01"
$ws.Range("E42").Value = "  -2.01%  "
$ws.Range("E43").Value = "  -2.51%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "126.44"
$ws.Range("E44").Value = "  -4.07%  "
$ws.Range("E45").Value = "  -2.13%  "
$ws.Range("E46").Value = "  -3.99%  "
$ws.Range("E47").Value = "  -3.01%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "17.08"
$ws.Range("E48").Value = "  -1.51%  "
$ws.Range("D49").Value = "1.732.58"
$ws.Range("E49").Value = "  -1.82%  "
$ws.Range("E50").Value = "  -1.53%  "
$ws.Range("E51").Value = "  -1.11%  "
